$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ------------------
# Overview sheet: zh-cn / de-de status cells (E2, F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) narrower
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn / de-de sheets: Status column (C) narrower
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
